$wb = $excel.ActiveWorkbook

function Set-CellAndHyperlink($ws, $addr, $val, $isHyperlink) {
    $ws.Range($addr).Value = $val
    if ($isHyperlink) {
        $target = ('$' + $addr.Substring(0,1) + '$' + $addr.Substring(1))
        foreach ($hl in $ws.Hyperlinks) {
            $a = $hl.Range.Address()
            if ($a -eq $target) {
                $hl.TextToDisplay = $val
            }
        }
    }
}


$ws = $wb.Worksheets.Item("Overview")
Set-CellAndHyperlink $ws "A2" "ffffe73ca54b-970e-4d2a-a723-8512024a563e.md" $true
Set-CellAndHyperlink $ws "D2" "2016-03-23 05:12:19" $false
Set-CellAndHyperlink $ws "A3" "ffffffd7a73dec-8c6f-4374-996c-a0ce6d4f22e1.md" $true
Set-CellAndHyperlink $ws "A4" "5bfe15cf-a495-4b59-84db-ca09775ae0ec.md" $true
Set-CellAndHyperlink $ws "B4" "Ready for handoff" $false
Set-CellAndHyperlink $ws "C4" "Ready for handoff" $false
Set-CellAndHyperlink $ws "D4" "2016-03-23 05:15:11" $false

$ws = $wb.Worksheets.Item("zh-cn")
Set-CellAndHyperlink $ws "A2" "ffffe73ca54b-970e-4d2a-a723-8512024a563e.md" $true
Set-CellAndHyperlink $ws "D2" "1e2b5805-6f60-4125-a897-ef3151d8ab4c.056ac7fcbe6e14b6529a7349561b36bd236bafa1.zh-cn.xlf" $true
Set-CellAndHyperlink $ws "E2" "2016-03-23 05:12:15" $false
Set-CellAndHyperlink $ws "F2" "1e2b5805-6f60-4125-a897-ef3151d8ab4c.md" $true
Set-CellAndHyperlink $ws "G2" "1e2b5805-6f60-4125-a897-ef3151d8ab4c.056ac7fcbe6e14b6529a7349561b36bd236bafa1.zh-cn.xlf" $true
Set-CellAndHyperlink $ws "H2" "2016-03-23 05:12:38" $false
Set-CellAndHyperlink $ws "A3" "ffffffd7a73dec-8c6f-4374-996c-a0ce6d4f22e1.md" $true
Set-CellAndHyperlink $ws "A4" "5bfe15cf-a495-4b59-84db-ca09775ae0ec.md" $true
Set-CellAndHyperlink $ws "C4" "Ready for handoff" $false
Set-CellAndHyperlink $ws "D4" "5bfe15cf-a495-4b59-84db-ca09775ae0ec.08e71ffab5dbe8be505740817c0a9cbcfc5cb687.zh-cn.xlf" $true
Set-CellAndHyperlink $ws "E4" "2016-03-23 05:15:05" $false
Set-CellAndHyperlink $ws "F4" "5bfe15cf-a495-4b59-84db-ca09775ae0ec.md" $true
Set-CellAndHyperlink $ws "G4" "5bfe15cf-a495-4b59-84db-ca09775ae0ec.08e71ffab5dbe8be505740817c0a9cbcfc5cb687.zh-cn.xlf" $true
Set-CellAndHyperlink $ws "H4" "2016-03-23 05:14:11" $false

$ws = $wb.Worksheets.Item("de-de")
Set-CellAndHyperlink $ws "A2" "ffffe73ca54b-970e-4d2a-a723-8512024a563e.md" $true
Set-CellAndHyperlink $ws "D2" "1e2b5805-6f60-4125-a897-ef3151d8ab4c.056ac7fcbe6e14b6529a7349561b36bd236bafa1.de-de.xlf" $true
Set-CellAndHyperlink $ws "E2" "2016-03-23 05:12:19" $false
Set-CellAndHyperlink $ws "F2" "1e2b5805-6f60-4125-a897-ef3151d8ab4c.md" $true
Set-CellAndHyperlink $ws "G2" "1e2b5805-6f60-4125-a897-ef3151d8ab4c.056ac7fcbe6e14b6529a7349561b36bd236bafa1.de-de.xlf" $true
Set-CellAndHyperlink $ws "H2" "2016-03-23 05:12:45" $false
Set-CellAndHyperlink $ws "A3" "ffffffd7a73dec-8c6f-4374-996c-a0ce6d4f22e1.md" $true
Set-CellAndHyperlink $ws "A4" "5bfe15cf-a495-4b59-84db-ca09775ae0ec.md" $true
Set-CellAndHyperlink $ws "C4" "Ready for handoff" $false
Set-CellAndHyperlink $ws "D4" "5bfe15cf-a495-4b59-84db-ca09775ae0ec.08e71ffab5dbe8be505740817c0a9cbcfc5cb687.de-de.xlf" $true
Set-CellAndHyperlink $ws "E4" "2016-03-23 05:15:11" $false
Set-CellAndHyperlink $ws "F4" "5bfe15cf-a495-4b59-84db-ca09775ae0ec.md" $true
Set-CellAndHyperlink $ws "G4" "5bfe15cf-a495-4b59-84db-ca09775ae0ec.08e71ffab5dbe8be505740817c0a9cbcfc5cb687.de-de.xlf" $true
Set-CellAndHyperlink $ws "H4" "2016-03-23 05:14:18" $false
